$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A23").Value = "edit1"
$ws.Range("B23").Value = "riya-morankar"
$ws.Range("C23").Value = "Merged"
$ws.Range("D23").Value = "N/A"

# "2025-06-20" looks like a date, and a plain .Value assignment would let
# Excel auto-convert it to a date serial. The other rows in this log store
# their Date column as literal text, so force the cell to Text format first
# (matches how the rest of column E was authored) before writing the value.
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "2025-06-20"

$ws.Range("F23").Value = "825d03f72949eefcc7953a7e836efe245eee87bb"
